{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Delete the empty paragraph and the \"Internal timer:\" paragraph that\n// followed \"From the class notes\" (the trailing two paragraphs of the doc).\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"\" || text === \"Internal timer:\") {\n    paragraphs.items[i].delete();\n  } else {\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"\" -or $t -eq \"Internal timer:\") {\n        $p.Range.Delete()\n    } else {\n        break\n    }\n}\n"}
